# "Put res and money options"
#
# This edit:
#   1. Adds two slide guides to the presentation (horizontal guide at 2160,
#      vertical guide at 2880 - PowerPoint's default guide color).
#   2. Retypes a few bullet paragraphs whose runs had been split up during
#      earlier authoring, so PowerPoint collapses them back into a single
#      run (the visible text does not change, only how it is split across
#      <a:r> runs).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide guides (View > Guides). Best effort: some hosts do not expose
#    a usable Guides collection; guard so the rest of the script always
#    runs even if this particular object model piece is unavailable.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($null -ne $guides) {
        $guides.Add(1, 2160)   # ppHorizontalGuide
        $guides.Add(2, 2880)   # ppVerticalGuide
    }
} catch {
}

# ---------------------------------------------------------------------
# Helper: collapse a paragraph's first $RunCount runs into a single run
# carrying $NewText (and the formatting of the paragraph's first run),
# exactly like re-typing the whole paragraph in the UI.
# ---------------------------------------------------------------------
function Merge-ParagraphRuns {
    param($Shape, $ParaIndex, $RunCount, $NewText)

    $tr = $Shape.TextFrame.TextRange
    $para = $tr.Paragraphs($ParaIndex, 1)

    $firstRun = $para.Runs(1, 1)
    $firstRun.Text = $NewText

    for ($i = 2; $i -le $RunCount; $i++) {
        $tr2 = $Shape.TextFrame.TextRange
        $para2 = $tr2.Paragraphs($ParaIndex, 1)
        $extraRun = $para2.Runs(2, 1)
        $extraRun.Text = ""
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2 ("Purpose") - second bullet paragraph: 3 runs -> 1 run.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
Merge-ParagraphRuns $shp2 2 3 "The game is supposed to be a fun and educational, to teach people about plant life in unnatural environment in an unique and fun way"

# ---------------------------------------------------------------------
# 3) Slide 4 ("Our plans for the future") - last bullet: 2 runs -> 1 run.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
Merge-ParagraphRuns $shp4 6 2 "Enable trading "

# ---------------------------------------------------------------------
# 4) Slide 5 ("What have we accomplished?") - second bullet, leading
#    "And " + "a cactus" runs -> single "And a cactus" run (the trailing
#    " " and ":" runs are left untouched).
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
Merge-ParagraphRuns $shp5 2 2 "And a cactus"
